$d = $word.ActiveDocument

# --- Paragraph 5: "Die Seite soll in einen öffentlichen, einen durch LogIn ..." ---
# Replace with plain text (drops the spell-check proofErr markers around LogIn/Login).
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range
[void]$r5.MoveEnd(1, -1)
[void]$r5.Delete()
$p5 = $d.Paragraphs(5)
[void]$p5.Range.InsertAfter("Die Seite soll in einen öffentlichen, einen durch Login geschützten und einen administrativen Bereich unterteilt sein, der nur per Login für registrierte Kunden erreichbar ist.")

# --- Paragraph 6: "Für den öffentlichen Bereich ..." ---
# Replace its text (removing the trailing "Außerdem ..." sentence), keep the
# trailing _GoBack bookmark in place at the end of the paragraph.
$p6 = $d.Paragraphs(6)
$r6 = $p6.Range
[void]$r6.MoveEnd(1, -1)
[void]$r6.Delete()
$p6 = $d.Paragraphs(6)
$r6b = $p6.Range
[void]$r6b.MoveEnd(1, -1)
[void]$r6b.InsertBefore("Für den öffentlichen Bereich sind beliebig viele Seiten vorgesehen, die allgemeine Informationen und Zugänge zu bestimmten Dokumenten (z.B. Handbüchern) zur HeatBox beinhalten. Der Administrator soll die Möglichkeit haben, die Homepage zu verändern um beispielsweise Seiten hinzuzufügen, um neue Versionen des Hardwareprodukts „HeatBox“ zu präsentieren. Die Seiten sollen in einer Hierarchie strukturiert werden und nach dieser in einem Navigationsbereich auch auswählbar sein. Die einzelnen Seiten sollen die Funktionalitäten erfüllen können, Titel bzw. Überschriften anzuzeigen, Texte darzustellen, Bilder anzuzeigen, eingebettete Videos wiederzugeben und Verknüpfungslinks zu beispielsweise Downloads bereitzustellen.")

# --- New paragraph 7: "Weiterhin soll im öffentlichen Bereich ..." ---
$p6 = $d.Paragraphs(6)
[void]$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs(7)
$p7.Range.Text = "Weiterhin soll im öffentlichen Bereich die Möglichkeit bestehen, auf möglichst unkomplizierte Weise den Admin zu kontaktieren, ein Impressum einzusehen und über einen Verknüpfungslink den Onlineshop des Kunden (externe Website) zu erreichen. "

# --- New paragraph 8: "Die Homepage soll ebenfalls, ..." ---
$p7 = $d.Paragraphs(7)
[void]$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)
$p8.Range.Text = "Die Homepage soll ebenfalls, gerecht der Anforderungen unserer Lehrerschaft, die Funktion besitzen, mithilfe eines Buttons die Anzeigesprachen Deutsch oder Englisch auszuwählen."

# --- New paragraph inserted between the two trailing empty paragraphs:
#     "Auf den einzelnen Seiten sollen dem Admin ..." ---
$pEmpty = $d.Paragraphs(9)
[void]$pEmpty.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs(10)
$pNew.Range.Text = "Auf den einzelnen Seiten sollen dem Admin Möglichkeiten bereitstehen einen Titel bzw. eine Überschrift festzulegen, einen Text zu verfassen, Bilder einzufügen, Videos einzubetten und Verknüpfungslinks zu beispielsweise Downloads anzugeben."
